$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8 from 45184 to 45185
$ws.Range("C2:C8").Value = 45185
